# Edit script: apply "Doing Updates for Financials" commit
# - Insert a new column D (shifting existing D:L data right by one)
# - Populate new column D with FY2018 data (period ending 2018-12-31, serial 43465)
# - Apply a handful of explicit data corrections to cells that don't follow
#   the simple shift pattern (re-stated prior-year figures)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before D; this shifts old D:K -> E:L
$ws.Columns("D:D").Insert()

# 2) The inserted column inherits formatting from column C by default.
#    Copy formats (number format/style) from column E (the old column D,
#    now shifted one to the right) back into the new column D so the new
#    column matches the rest of the data columns (date format on row 7/38/80,
#    number format elsewhere).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Populate the new column D with the FY2018 values for each data row.
$newColD = @{
    7 = 43465
    8 = 1708900
    9 = "NA"
    10 = "NA"
    12 = "NA"
    13 = 0
    14 = 0
    15 = 0
    17 = 572500
    18 = 1136500
    20 = -463000
    21 = 687500
    22 = 0
    23 = 673500
    24 = 168100
    25 = 0
    26 = 505300
    27 = 504400
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 463000
    33 = 504400
    34 = 0
    35 = 504400
    38 = 43465
    41 = 269200
    42 = 312900
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 59100
    49 = 0
    50 = 0
    51 = 0
    52 = 9400
    53 = 0
    54 = 47364800
    57 = 0
    58 = 0
    59 = 530700
    60 = 0
    61 = 258200
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 42957700
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 2730900
    73 = 0
    74 = 0
    75 = 0
    76 = 4407100
    77 = 0
    80 = 43465
    81 = 504400
    83 = 14000
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 651600
    91 = -11500
    92 = 0
    93 = 0
    94 = -4289100
    96 = -62000
    97 = 0
    98 = 0
    99 = 0
    100 = 3619400
    101 = 0
    102 = -18200
}

foreach ($row in $newColD.Keys) {
    $cell = $ws.Cells.Item([int]$row, 4)
    $cell.Value2 = $newColD[$row]
}

# 4) A few cells were re-stated as part of this update and do not match a
#    pure column shift of the previous values; set those explicitly.
$corrections = @(
    @{ Row = 24; Col = "E"; Val = 190000 }
    @{ Row = 26; Col = "E"; Val = 385300 }
    @{ Row = 27; Col = "E"; Val = 385300 }
    @{ Row = 89; Col = "E"; Val = 546000 }
    @{ Row = 89; Col = "F"; Val = 565000 }
    @{ Row = 94; Col = "E"; Val = -4341300 }
    @{ Row = 94; Col = "F"; Val = -5418400 }
    @{ Row = 100; Col = "E"; Val = 3591800 }
    @{ Row = 100; Col = "F"; Val = 5050800 }
)

foreach ($fix in $corrections) {
    $ws.Range("$($fix.Col)$($fix.Row)").Value2 = $fix.Val
}

$wb.Save()
